$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mae")

# Update the descriptive text in the shared string table (cell O3)
$ws.Range("O3").Value = "Анализ признаков в датасете (размер выборки 1925 примеров)"

# Update numeric values for P and S columns, rows 5-9
$ws.Range("P5").Value = 601.5
$ws.Range("S5").Value = 888

$ws.Range("P6").Value = 38.01
$ws.Range("S6").Value = 4

$ws.Range("P7").Value = 45.61
$ws.Range("S7").Value = 762

$ws.Range("P8").Value = 0.59
$ws.Range("S8").Value = 1229

$ws.Range("P9").Value = 0.2
$ws.Range("S9").Value = 1597

# Only S10 changes in row 10 (P10 stays the same)
$ws.Range("S10").Value = 1697

# Update the active selection shown when the sheet was last saved
$ws.Range("N9").Select()
